$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Text relabelling (shared text used in multiple cells) ---
$ws.Range("C16").Value = "Semakan Kali Pertama"
$ws.Range("C22").Value = "Semakan Kali Pertama"
$ws.Range("C28").Value = "Semakan Kali Pertama"

$ws.Range("C17").Value = "Semakan Kali Kedua"
$ws.Range("C23").Value = "Semakan Kali Kedua"
$ws.Range("C29").Value = "Semakan Kali Kedua"

$ws.Range("C18").Value = "Semakan Kali Ketiga"
$ws.Range("C24").Value = "Semakan Kali Ketiga"
$ws.Range("C30").Value = "Semakan Kali Ketiga"

$ws.Range("C19").Value = "Semakan Kali Keempat"
$ws.Range("C25").Value = "Semakan Kali Keempat"
$ws.Range("C31").Value = "Semakan Kali Keempat"

$ws.Range("C34").Value = "Bouquet Kreatif"
$ws.Range("C35").Value = "Kad Raya Untuk Guruku"
$ws.Range("C36").Value = "Riang Ria Kuih Raya"
$ws.Range("C37").Value = "Creative Collage"

# --- New competition rows (previously blank) ---
$ws.Range("C38").Value = "Bowling Padang"
$ws.Range("C39").Value = "Theme Party"
$ws.Range("C40").Value = "Melukis Poster Koop"

# --- Updated figures / new data ---
$ws.Range("D18").Value = 11110
$ws.Range("E18").Value = 150
$ws.Range("E23").Value = 400
$ws.Range("E29").Value = 1500
$ws.Range("D38").Value = 100
$ws.Range("D39").Value = 100
$ws.Range("D40").Value = 200

# --- Move statement title from E4 to D4 and merge D4:G4 ---
$ws.Range("D4").Value = $ws.Range("E4").Value
$ws.Range("E4").ClearContents()
$ws.Range("D4:G4").Merge()

# --- Merge-range adjustments ---
$ws.Range("B15:C15").UnMerge()
$ws.Range("B12:F12").Merge()
$ws.Range("B21:E21").Merge()
$ws.Range("B27:E27").Merge()
$ws.Range("B33:E33").Merge()
$ws.Range("B43:E43").Merge()

# --- Remove trailing filler row ---
$ws.Rows.Item(1001).Delete()

# --- Page setup / print options ---
$ws.PageSetup.Zoom = $false
$ws.PageSetup.FitToPagesWide = 1
$ws.PageSetup.FitToPagesTall = 1
$ws.PageSetup.HeaderMargin = 0
$ws.PageSetup.FooterMargin = 0
$ws.PageSetup.CenterHorizontally = $true

# --- Picture reposition/resize (col B, small inset near top-left) ---
$shp = $ws.Shapes.Item(1)
$shp.Left = 41.2125
$shp.Top = 14.25
$shp.Width = 46.5
$shp.Height = 47.25

# --- Reset selection to default (A1) ---
$ws.Range("A1").Select()
